# CrewAI Robust Backend Ready!
# Relabel headers, reshuffle the footprint columns (D/E/F), and annotate
# column headers with data-type comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename / lower-case the header row -------------------------------
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- 2. Shift the D/E/F data columns and populate the new F values -------
# D <- old E, E <- old F, F <- new "climate change" figures
# (index 0 corresponds to row 2, index 6 to row 8)
$climateChange = @(0.0000025775142, 0.0000009455103799999999, 0.0000023034924, 0.0000032414908, 0.0000029040017, 0.0000028833396, 0.0000022716888)

for ($row = 2; $row -le 8; $row++) {
    $oldE = $ws.Cells.Item($row, 5).Value2
    $oldF = $ws.Cells.Item($row, 6).Value2

    $ws.Cells.Item($row, 4).Value = $oldE
    $ws.Cells.Item($row, 5).Value = $oldF
    $ws.Cells.Item($row, 6).Value = $climateChange[$row - 2]
}

# --- 3. Add header comments describing the data type ---------------------
$comments = @{
    "A1" = "Data type: Categorical (text)"
    "B1" = "Data type: Various (e.g. kg, kWh)"
    "C1" = "Data type: Categorical (text)"
    "D1" = "Data type: Carbon footprint"
    "E1" = "Data type: Cumulative energy demand"
    "F1" = "Data type: Climate change impact"
    "G1" = "Data type: Categorical (text)"
}

foreach ($ref in @("A1", "B1", "C1", "D1", "E1", "F1", "G1")) {
    $ws.Range($ref).AddComment($comments[$ref]) | Out-Null
}
